$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update legend cells in row 1 (B1 and D1 keep their text; C1's legend
# label changes from "ordered-DenStream" to "DistStream-DenStream")
$ws.Range("B1").Value = "MOA=1.0"
$ws.Range("C1").Value = "DistStream-DenStream"
$ws.Range("D1").Value = "unordered-DenStream"

# Update column widths: B -> ~12.33203125, C -> 21, D stays 29.1640625
$ws.Columns.Item(2).ColumnWidth = 11.571428571428571
$ws.Columns.Item(3).ColumnWidth = 20.285714285714285

# Update active selection to D3
$ws.Range("D3").Select()
